$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.034.56"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.910.53"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'0.7854"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.13%  "
$ws.Range("D6").Value = "'241.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.3156"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "'26.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").Value = "'0.06911"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "'0.07968"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "1.908.78"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "'0.7439"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "'5.216"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "'93.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "30.046.86"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'13.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'5.887"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.83%  "
$ws.Range("D19").Value = "'246.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("D20").Value = "'0.000007768"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "2.150.29"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'6.879"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").Value = "'169.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'9.283"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "'0.1376"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.57%  "
$ws.Range("D28").Value = "'18.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "'2.033"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "'1.376"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").Value = "'4.330"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "'4.087"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "'0.05467"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").Value = "'1.258"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "'0.01937"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "'2.791"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").Value = "'6.143"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "'0.4426"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "'0.8365"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "'7.532"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").Value = "'9.756"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "'981.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.54%  "
$ws.Range("D50").Value = "2.057.80"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'36.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
